$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Internet"/G column entirely, shifting everything left
$ws.Range("G:G").Delete()

# Update header row (B1:F1) with new keyword set
$ws.Range("B1").Value = "Processing"
$ws.Range("C1").Value = "AI"
$ws.Range("D1").Value = "Technology"
$ws.Range("E1").Value = "Data"
$ws.Range("F1").Value = "Adjectives"

# Update the data rows (B2:F7) with the new counts
$data = @(
    @(0, 181, 76, 0, 0),
    @(0, 186, 75, 0, 2),
    @(0, 303, 103, 29, 4),
    @(0, 266, 105, 19, 4),
    @(0, 245, 100, 22, 4),
    @(0, 323, 132, 28, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
